$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maria Clara's availability flips from FALSE to TRUE
$ws.Range("C3").Value = $true

# Claudio Golveia - new row 8
$a8 = $ws.Range("A8")
$a8.Value = "'154"
$a8.ClearFormats()
$ws.Range("B8").Value = "Claudio Golveia"
$ws.Range("C8").Value = $true

# Roberta Miranda - new row 9
$a9 = $ws.Range("A9")
$a9.Value = "'15"
$a9.ClearFormats()
$ws.Range("B9").Value = "Roberta Miranda"
$ws.Range("C9").Value = $true
